$wb = $excel.ActiveWorkbook

# --- SCHEME_MEASURES sheet: renumber indicators MQMS0x -> MQME00x ---
$wsMeasures = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsMeasures.Range("A2").Value = "MQME001"
$wsMeasures.Range("A3").Value = "MQME002"
$wsMeasures.Range("A4").Value = "MQME003"
$wsMeasures.Range("A5").Value = "MQME004"
$wsMeasures.Range("A6").Value = "MQME005"

# --- METADATA_ISSUES sheet: renumber indicator MQME01 -> MQME008 ---
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
$wsIssues.Range("A2").Value = "MQME008"

# --- METADATA_MEASURES sheet: renumber + relabel, drop last row ---
$wsMetaMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMetaMeasures.Range("A2").Value = "MQME006"
$wsMetaMeasures.Range("B2").Value = "Total number of length-required columns"
$wsMetaMeasures.Range("C2").Value = 28

$wsMetaMeasures.Range("A3").Value = "MQME007"
$wsMetaMeasures.Range("B3").Value = "Total number of NUMBER columns"
$wsMetaMeasures.Range("C3").Value = 218

$wsMetaMeasures.Range("A4:C4").Delete()

# --- METADATA_METRICS sheet: renumber + relabel + add new indicator rows ---
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")

# Percentage-looking values must stay plain text (as in the source file),
# so force text formatting on column C before writing them - otherwise
# Excel auto-converts "100.00%" into the number 1 with a percent format.
$wsMetrics.Range("C2:C12").NumberFormat = "@"

$wsMetrics.Range("A2").Value = "MQID001"
$wsMetrics.Range("B2").Value = "Table names in singular"
$wsMetrics.Range("C2").Value = "100.00%"

$wsMetrics.Range("A3").Value = "MQID002"
$wsMetrics.Range("B3").Value = "Table with recommended name length"
$wsMetrics.Range("C3").Value = "100.00%"

$wsMetrics.Range("A4").Value = "MQID003"
$wsMetrics.Range("B4").Value = "Columns with correct prefixes"
$wsMetrics.Range("C4").Value = "100.00%"

$wsMetrics.Range("A5").Value = "MQID004"
$wsMetrics.Range("B5").Value = "Columns with recommended name size"
$wsMetrics.Range("C5").Value = "100.00%"

$wsMetrics.Range("A6").Value = "MQID005"
$wsMetrics.Range("B6").Value = "Columns with comments"
$wsMetrics.Range("C6").Value = "99.59%"

$wsMetrics.Range("A7").Value = "MQID006"
$wsMetrics.Range("B7").Value = "Table with standard PK prefixes"
$wsMetrics.Range("C7").Value = "100.00%"

$wsMetrics.Range("A8").Value = "MQID007"
$wsMetrics.Range("B8").Value = "Table with standard FK prefixes"
$wsMetrics.Range("C8").Value = "0.00%"

$wsMetrics.Range("A9").Value = "MQID008"
$wsMetrics.Range("B9").Value = "Table with standard UK prefixes"
$wsMetrics.Range("C9").Value = "0.00%"

$wsMetrics.Range("A10").Value = "MQID009"
$wsMetrics.Range("B10").Value = "NUMBER columns with valid scale"
$wsMetrics.Range("C10").Value = "100.00%"

$wsMetrics.Range("A11").Value = "MQID010"
$wsMetrics.Range("B11").Value = "Columns with valid num_distinct"
$wsMetrics.Range("C11").Value = "100.00%"

$wsMetrics.Range("A12").Value = "MQID011"
$wsMetrics.Range("B12").Value = "Columns with valid num_nulls"
$wsMetrics.Range("C12").Value = "100.00%"
